# The workbook had an empty row 2 between the header row (row 1) and the
# start of the DAX reference table (which began at row 3). This edit
# removes that stray blank row so the table starts immediately below the
# header, shifting all subsequent rows up by one (old row 3 -> row 2,
# old row 22 -> row 21, dimension A1:B22 -> A1:B21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2").Delete()

# Reflect the author's final selection/cursor position in the saved file.
$ws.Range("A2").Select() | Out-Null
